$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell reference -> new text value.
# Values that look like plain numbers are given a leading apostrophe
# so Excel stores them as text (matching the original inlineStr cells)
# instead of auto-converting them to numeric cells.
$updates = @{
    'D2' = '30.203.79'
    'E2' = '  -0.65%  '
    'D3' = '1.851.81'
    'E3' = '  -2.10%  '
    'D4' = '''0.9995'
    'E4' = '  -0.18%  '
    'D5' = '''236.22'
    'E5' = '  -0.65%  '
    'D6' = '''0.9989'
    'E6' = '  -0.16%  '
    'D7' = '''0.4783'
    'E7' = '  -2.30%  '
    'D8' = '''0.2809'
    'D9' = '''0.06482'
    'E9' = '  -3.19%  '
    'D10' = '1.855.62'
    'E10' = '  -2.00%  '
    'D11' = '''0.07311'
    'E11' = '  -0.49%  '
    'E12' = '  -3.87%  '
    'D13' = '''5.114'
    'E13' = '  -0.09%  '
    'D14' = '''87.23'
    'E14' = '  -0.33%  '
    'D15' = '''0.6469'
    'E15' = '  -2.37%  '
    'D16' = '30.148.56'
    'E16' = '  -0.80%  '
    'E17' = '  -1.32%  '
    'D18' = '''0.9990'
    'E18' = '  -0.22%  '
    'D19' = '''0.000007623'
    'E19' = '  -2.52%  '
    'D20' = '''225.17'
    'E20' = '  +18.78%  '
    'D21' = '2.099.56'
    'E21' = '  -1.43%  '
    'D22' = '''5.286'
    'E22' = '  -0.47%  '
    'D23' = '''0.9999'
    'E23' = '  -0.13%  '
    'D24' = '''6.079'
    'E24' = '  -0.37%  '
    'D25' = '''9.225'
    'E25' = '  -2.47%  '
    'D26' = '''163.58'
    'E26' = '  +0.17%  '
    'D27' = '''18.54'
    'E27' = '  +1.69%  '
    'D28' = '''1.916'
    'E28' = '  -0.55%  '
    'D29' = '''1.430'
    'E29' = '  -2.37%  '
    'B30' = 'InternetComputer(DFINITY)'
    'C30' = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
    'D30' = '''4.247'
    'E30' = '  -2.47%  '
    'B31' = 'Stellar'
    'C31' = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
    'D31' = '''0.09193'
    'E31' = '  +0.52%  '
    'D32' = '''3.955'
    'E32' = '  -1.91%  '
    'E33' = '  -3.48%  '
    'D34' = '''0.7396'
    'E34' = '  +0.02%  '
    'D35' = '''1.143'
    'E35' = '  +4.18%  '
    'D36' = '''2.685'
    'E36' = '  -1.14%  '
    'D37' = '''0.01812'
    'E37' = '  +0.16%  '
    'D38' = '''2.605'
    'E38' = '  -2.30%  '
    'D39' = '''0.9070'
    'E39' = '  -1.54%  '
    'D40' = '''2.052'
    'E40' = '  +0.88%  '
    'D41' = '''5.953'
    'E41' = '  +0.54%  '
    'D42' = '''106.43'
    'E42' = '  +0.43%  '
    'D43' = '''0.4254'
    'E43' = '  -2.94%  '
    'D44' = '''0.9986'
    'E44' = '  +0.67%  '
    'B45' = 'Algorand'
    'C45' = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
    'D45' = '''0.1321'
    'E45' = '  -3.51%  '
    'B46' = 'Aptos'
    'C46' = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
    'D46' = '''7.382'
    'E46' = '  -2.39%  '
    'D47' = '''1.551'
    'E47' = '  +11.05%  '
    'E48' = '  -6.36%  '
    'D49' = '''34.22'
    'E49' = '  -1.76%  '
    'D50' = '''8.758'
    'E50' = '  -2.54%  '
    'D51' = '''0.05655'
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
